$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-* annotations are re-curated from "dimension" to "measure"
$ws.Range("D2").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("E2").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("J2").Value = "iaest-measure:sexo"
$ws.Range("K2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("M2").Value = "iaest-measure:nacionalidad-area-nombre"

# Row 3: "dim" role becomes "medida" for the same columns
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("J3").Value = "medida"
$ws.Range("K3").Value = "medida"
$ws.Range("M3").Value = "medida"

# Row 4: datatype becomes "xsd:int" instead of URI-* / skos:Concept
$ws.Range("D4").Value = "xsd:int"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("M4").Value = "xsd:int"

# Row 5: mapping files for these columns no longer apply - remove the cells entirely
$ws.Range("E5").Clear()
$ws.Range("J5").Clear()
$ws.Range("M5").Clear()
